# Add team record (Wins/Losses/Ties) columns to the SFG 1990 roster sheet.
# New columns AC, AD, AE are appended after the existing data (A1:AB53),
# with header labels in row 1 and the team's 85-77-0 record repeated on
# every player row (2-53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the existing header formatting (bold,
# bordered, centered) from the last header cell (AB1) onto the three new
# header cells, then set their labels.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# --- Data rows (2-53): same team record (85 wins, 77 losses, 0 ties) for
# every player row.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 29).Value = 85   # column AC - Wins
    $ws.Cells.Item($r, 30).Value = 77   # column AD - Losses
    $ws.Cells.Item($r, 31).Value = 0    # column AE - Ties
}
